$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "time_taken" header in F1, copying the header formatting from
# E1 (bold font, border, centered) so F1 reuses the existing header style
# instead of Excel minting a brand-new style entry.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1, 6).Value = "time_taken"

# time_taken values for data rows 2-66, in row order
$timeTakenValues = @(
    "2021-10-05 13:42:25.371973"
    "2021-10-05 13:42:25.371986"
    "2021-10-05 13:42:25.371990"
    "2021-10-05 13:42:25.371993"
    "2021-10-05 13:42:25.371997"
    "2021-10-05 13:42:25.372000"
    "2021-10-05 13:42:25.372003"
    "2021-10-05 13:42:25.372006"
    "2021-10-05 13:42:25.372009"
    "2021-10-05 13:42:25.372012"
    "2021-10-05 13:42:25.372015"
    "2021-10-05 13:42:25.372018"
    "2021-10-05 13:42:25.372021"
    "2021-10-05 13:42:25.372024"
    "2021-10-05 13:42:25.372027"
    "2021-10-05 13:42:25.372030"
    "2021-10-05 13:42:25.372033"
    "2021-10-05 13:42:25.372036"
    "2021-10-05 13:42:25.372039"
    "2021-10-05 13:42:25.372042"
    "2021-10-05 13:42:25.372045"
    "2021-10-05 13:42:25.372048"
    "2021-10-05 13:42:25.372051"
    "2021-10-05 13:42:25.372054"
    "2021-10-05 13:42:25.372057"
    "2021-10-05 13:42:25.372061"
    "2021-10-05 13:42:25.372063"
    "2021-10-05 13:42:25.372067"
    "2021-10-05 13:42:25.372069"
    "2021-10-05 13:42:25.372072"
    "2021-10-05 13:42:25.372075"
    "2021-10-05 13:42:25.372078"
    "2021-10-05 13:42:25.372082"
    "2021-10-05 13:42:25.372085"
    "2021-10-05 13:42:25.372088"
    "2021-10-05 13:42:25.372091"
    "2021-10-05 13:42:25.372094"
    "2021-10-05 13:42:25.372097"
    "2021-10-05 13:42:25.372100"
    "2021-10-05 13:42:25.372103"
    "2021-10-05 13:42:25.372107"
    "2021-10-05 13:42:25.372110"
    "2021-10-05 13:42:25.372113"
    "2021-10-05 13:42:25.372116"
    "2021-10-05 13:42:25.372119"
    "2021-10-05 13:42:25.372121"
    "2021-10-05 13:42:25.372124"
    "2021-10-05 13:42:25.372127"
    "2021-10-05 13:42:25.372131"
    "2021-10-05 13:42:25.372134"
    "2021-10-05 13:42:25.372137"
    "2021-10-05 13:42:25.372139"
    "2021-10-05 13:42:25.372143"
    "2021-10-05 13:42:25.372148"
    "2021-10-05 13:42:25.372151"
    "2021-10-05 13:42:25.372154"
    "2021-10-05 13:42:25.372157"
    "2021-10-05 13:42:25.372160"
    "2021-10-05 13:42:25.372163"
    "2021-10-05 13:42:25.372166"
    "2021-10-05 13:42:25.372169"
    "2021-10-05 13:42:25.372172"
    "2021-10-05 13:42:25.372175"
    "2021-10-05 13:42:25.372178"
    "2021-10-05 13:42:25.372183"
)

$startRow = 2
for ($i = 0; $i -lt $timeTakenValues.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 6).Value = $timeTakenValues[$i]
}

Write-Host "Updated F1:F66 with time_taken column"
